$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.685.29"
$ws.Range("E2").Value = "  +1.32%  "
$ws.Range("D3").Value = "1.632.18"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("D5").Value = "213.42"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("E6").Value = "  +3.49%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  +1.49%  "
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("D10").Value = "19.25"
$ws.Range("E10").Value = "  +2.83%  "
$ws.Range("D11").Value = "0.0840"
$ws.Range("E11").Value = "  +3.26%  "
$ws.Range("D12").Value = "1.859.28"
$ws.Range("D13").Value = "1.642.27"
$ws.Range("E13").Value = "  +1.25%  "
$ws.Range("E14").Value = "  +1.34%  "
$ws.Range("E15").Value = "  +0.97%  "
$ws.Range("D16").Value = "26.683.03"
$ws.Range("E16").Value = "  +1.29%  "
$ws.Range("D17").Value = "63.48"
$ws.Range("E17").Value = "  +1.61%  "
$ws.Range("D18").Value = "0.0₃0741"
$ws.Range("E18").Value = "  +2.12%  "
$ws.Range("D19").Value = "218.92"
$ws.Range("E19").Value = "  +8.39%  "
$ws.Range("E21").Value = "  +0.98%  "
$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").Value = "9.37"
$ws.Range("E22").Value = "  +0.82%  "
$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").Value = "6.16"
$ws.Range("E23").Value = "  +1.70%  "
$ws.Range("D24").Value = "1.96"
$ws.Range("E24").Value = "  +4.60%  "
$ws.Range("D25").Value = "147.66"
$ws.Range("E25").Value = "  +2.35%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  +1.58%  "
$ws.Range("D28").Value = "6.86"
$ws.Range("E28").Value = "  +4.36%  "
$ws.Range("D29").Value = "15.56"
$ws.Range("E29").Value = "  +2.38%  "
$ws.Range("E30").Value = "  -2.56%  "
$ws.Range("E31").Value = "  +0.20%  "
$ws.Range("D32").Value = "3.31"
$ws.Range("E32").Value = "  +4.03%  "
$ws.Range("D33").Value = "2.99"
$ws.Range("E33").Value = "  +2.41%  "
$ws.Range("D34").Value = "1.51"
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("E35").Value = "  +0.29%  "
$ws.Range("D36").Value = "1.216.29"
$ws.Range("E36").Value = "  +4.76%  "
$ws.Range("E37").Value = "  +4.31%  "
$ws.Range("E38").Value = "  +0.25%  "
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("E40").Value = "  +1.10%  "
$ws.Range("E41").Value = "  -2.02%  "
$ws.Range("D42").Value = "0.794"
$ws.Range("E42").Value = "  +1.50%  "
$ws.Range("E43").Value = "  -0.91%  "
$ws.Range("D44").Value = "1.771.18"
$ws.Range("E44").Value = "  +0.68%  "
$ws.Range("D45").Value = "92.85"
$ws.Range("E45").Value = "  +0.59%  "
$ws.Range("E46").Value = "  +2.67%  "
$ws.Range("D47").Value = "0.0₆0105"
$ws.Range("E47").Value = "  +0.56%  "
$ws.Range("D48").Value = "55.04"
$ws.Range("E48").Value = "  +2.29%  "
$ws.Range("E49").Value = "  +0.63%  "
$ws.Range("D50").Value = "7.61"
$ws.Range("D51").Value = "0.408"
$ws.Range("E51").Value = "  -0.44%  "
